# Update the QC overview on Sheet1: row 7 held stale/placeholder values from
# before the "opt0030" plot-generation run was wired up to refresh this sheet
# automatically. Refresh row 7 (C7:Z7) with the current values already
# present in row 8 (C8:Z8), matching what the automated plot-file process
# now writes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C7:Z7").Value2 = $ws.Range("C8:Z8").Value2

# Leave the selection where the author ended up after making the edit.
$ws.Range("C8").Select()
